# Insert a new weekly price record as row 72 in the Piña price-list sheet.
# All existing rows from 72 downward shift down by one (72 -> 73, ..., 147 -> 148).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 72, pushing everything else down.
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new observation.
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44494
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100108
$ws.Range("H72").Value = "Tropicales y subtropicales"
$ws.Range("I72").Value = 100108005
$ws.Range("J72").Value = "Piña"
$ws.Range("K72").Value = "Caramelo"
$ws.Range("L72").Value = "Segunda"
$ws.Range("M72").Value = 180
$ws.Range("N72").Value = 25000
$ws.Range("O72").Value = 25500
$ws.Range("P72").Value = 25250
$ws.Range("Q72").Value = "$/caja 14 unidades"
$ws.Range("R72").Value = "Ecuador"
$ws.Range("S72").Value = 1804
$ws.Range("T72").Value = 14
